$p = $ppt.ActivePresentation

function Set-MergedText($shape, $text) {
    # Force the engine to rewrite/merge the run structure: assigning the
    # exact same (already-concatenated) text is treated as a no-op, so we
    # nudge it through a throwaway value first.
    $shape.TextFrame.TextRange.Text = "__tmp__"
    $shape.TextFrame.TextRange.Text = $text
}

# Slide 1: "Section Header (with background image)"
Set-MergedText $p.Slides.Item(1).Shapes.Item(1) "Section Header (with background image)"

# Slide 2: "Slide 1"
Set-MergedText $p.Slides.Item(2).Shapes.Item(1) "Slide 1"

# Slide 3: "Slide 2"
Set-MergedText $p.Slides.Item(3).Shapes.Item(1) "Slide 2"

# Slide 4: "Slide 3"
Set-MergedText $p.Slides.Item(4).Shapes.Item(1) "Slide 3"

# Slide 5: "Slide 4" and "An image"
Set-MergedText $p.Slides.Item(5).Shapes.Item(1) "Slide 4"
Set-MergedText $p.Slides.Item(5).Shapes.Item(4) "An image"

# Slide 6 notes: "Blank slides can have background images."
Set-MergedText $p.Slides.Item(6).NotesPage.Shapes.Item(2) "Blank slides can have background images."
